$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) cells that look numeric keep their exact
# original text representation (e.g. trailing zeros) by formatting
# the cells as Text before assigning the new value - mirrors the
# source data, which stores every Price/Volume cell as text.

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '54.151.44'
$ws.Range("E2").Value = '  -3.77%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.259.91'
$ws.Range("E3").Value = '  -4.83%  '

# Row 4
$ws.Range("E4").Value = '  -0.01%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '490.02'
$ws.Range("E5").Value = '  -2.23%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '126.92'
$ws.Range("E6").Value = '  -2.55%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.00'
$ws.Range("E7").Value = '  +0.06%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.521'
$ws.Range("E8").Value = '  -4.17%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.264.64'
$ws.Range("E9").Value = '  -4.95%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0924'
$ws.Range("E10").Value = '  -6.08%  '

# Row 11
$ws.Range("E11").Value = '  -1.34%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '4.73'
$ws.Range("E12").Value = '  +1.61%  '

# Row 13
$ws.Range("E13").Value = '  -2.33%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '2.661.70'
$ws.Range("E14").Value = '  -4.84%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '21.28'
$ws.Range("E15").Value = '  -1.45%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '54.167.28'
$ws.Range("E16").Value = '  -3.67%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.0000128'
$ws.Range("E17").Value = '  -2.60%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.249.43'
$ws.Range("E18").Value = '  -6.60%  '

# Row 19
$ws.Range("B19").Value = 'Polkadot'
$ws.Range("C19").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.99'
$ws.Range("E19").Value = '  -0.71%  '

# Row 20
$ws.Range("B20").Value = 'Chainlink'
$ws.Range("C20").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '9.63'
$ws.Range("E20").Value = '  -4.24%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '301.17'
$ws.Range("E21").Value = '  -2.04%  '

# Row 22
$ws.Range("E22").Value = '  -2.08%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.999'
$ws.Range("E23").Value = '  -0.20%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '63.89'
$ws.Range("E24").Value = '  -1.29%  '

# Row 25
$ws.Range("E25").Value = '  +0.17%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.368'
$ws.Range("E26").Value = '  -0.37%  '

# Row 27
$ws.Range("E27").Value = '  -2.43%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.09'
$ws.Range("E28").Value = '  -2.69%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '170.36'
$ws.Range("E29").Value = '  -1.06%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0₃0694'
$ws.Range("E30").Value = '  -2.82%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.60'
$ws.Range("E31").Value = '  -1.84%  '

# Row 32
$ws.Range("E32").Value = '  -0.12%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.76'
$ws.Range("E33").Value = '  +0.24%  '

# Row 34
$ws.Range("E34").Value = '  -0.07%  '

# Row 35
$ws.Range("E35").Value = '  -1.46%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '17.49'
$ws.Range("E36").Value = '  -0.58%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.16'
$ws.Range("E37").Value = '  -0.60%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.851'
$ws.Range("E38").Value = '  +7.35%  '

# Row 39
$ws.Range("E39").Value = '  -4.77%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '35.79'
$ws.Range("E40").Value = '  -0.64%  '

# Row 41
$ws.Range("B41").Value = 'PolygonEcosystemToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.369'
$ws.Range("E41").Value = '  -0.09%  '

# Row 42
$ws.Range("B42").Value = 'Stacks'
$ws.Range("C42").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.38'
$ws.Range("E42").Value = '  -2.44%  '

# Row 43
$ws.Range("E43").Value = '  -0.66%  '

# Row 44
$ws.Range("B44").Value = 'Aave'
$ws.Range("C44").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '122.35'
$ws.Range("E44").Value = '  -6.67%  '

# Row 45
$ws.Range("B45").Value = 'RenderToken'
$ws.Range("C45").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '4.70'
$ws.Range("E45").Value = '  -1.12%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0881'
$ws.Range("E46").Value = '  -2.43%  '

# Row 47
$ws.Range("E47").Value = '  -4.54%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '238.00'
$ws.Range("E48").Value = '  -1.62%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0473'
$ws.Range("E49").Value = '  -2.00%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0203'
$ws.Range("E50").Value = '  -2.83%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '16.27'
$ws.Range("E51").Value = '  -3.52%  '
